$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 3-14 (IP_NAS_AP, Inicio_de_Conexion_Dia, FIN_de_Conexion_Dia, Input_Octects, Output_Octects)
# Values are stored as text (leading apostrophe) to match the original all-text data layout
# and avoid Excel auto-converting date-like / numeric-like strings into dates or numbers.
$data = @(
    @("192.168.247.18","2019-08-21","2020-02-19","0","0"),
    @("192.168.247.12","2019-02-26","2019-02-26","3084","1344"),
    @("192.168.247.18","2019-04-03","2020-02-05","0","0"),
    @("192.168.247.11","2019-02-11","2019-03-14","1170","495"),
    @("192.168.247.20","2019-06-27","2020-10-08","8990","14372"),
    @("192.168.247.19","2019-02-19","2019-02-19","2554","334"),
    @("192.168.247.11","2019-02-12","2019-03-06","544518","1881843"),
    @("192.168.247.19","2019-04-03","2020-10-13","0","0"),
    @("192.168.247.18","2019-02-27","2020-02-04","2754","2432"),
    @("192.168.247.15","2019-03-28","2019-03-28","16893","23855"),
    @("192.168.247.16","2019-03-15","2019-05-02","0","0"),
    @("192.168.247.18","2019-06-24","2019-06-24","0","0")
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value = "'" + $row[4]
    $r = $r + 1
}
